# EDA and preprocessing code update
# Refresh the doc2vec + logistic-regression row of the statistics table
# (precision/recall/f1 per class + overall accuracy) after rerunning the
# preprocessing / EDA pipeline, and leave the selection where the author
# left off working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Row 3: A3="doc2vec", B3="logistic regression"
#   C3/D3 = precision (sarcasm / not sarcasm)
#   E3/F3 = recall (sarcasm / not sarcasm)
#   G3/H3 = f1-score (sarcasm / not sarcasm)
#   I3    = accuracy
$ws.Range("C3").Value = 0.61
$ws.Range("D3").Value = 0.62
$ws.Range("E3").Value = 0.63
$ws.Range("F3").Value = 0.6
$ws.Range("G3").Value = 0.62
$ws.Range("H3").Value = 0.61
$ws.Range("I3").Value = 0.61

# Move the active selection to match where the author left off
$ws.Range("H11").Select()
